$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D8").Value = "카카오브레인"

$ws.Range("D28").Value = "[임피던스 제어] Interaction Control 적용(1) -  Virtual Trajectory & nodic impedance"
$ws.Range("E28").Value = "https://ropiens.tistory.com/120"

$ws.Range("D46").Value = "[Bioinformatics] 2021년 05월, 유전체 빅데이터 분석 교육-""예비전문가 과정"" (9기) 21년도 교육생 모집 [한국바이오협회]"
$ws.Range("E46").Value = "https://bioinformaticsandme.tistory.com/395"

$ws.Range("D51").Value = "MySQL workbench에서 select로 조회했는데 row들이 잘 안 보일 때"
$ws.Range("E51").Value = "https://bskyvision.com/1193"
